$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 2.045107565409183 * [Math]::Pow(10, 21)
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 2.045107565409197 * [Math]::Pow(10, 21)
